$wb = $excel.ActiveWorkbook

# --- Re-add the "Scaled Beta" distribution row to the validation table ---
$wsVal = $wb.Worksheets.Item("Internal - Data Validation")

# Insert a new row at row 6, pushing the existing rows (Log Uniform, etc.) down.
$wsVal.Rows.Item(6).Insert()

# Pick up the same cell formatting used by the other 4-parameter rows (e.g. row 4,
# "Truncated Normal": Min(a)/Max(b)/<p3>/<p4> with G:H left blank) for the new row.
$wsVal.Range("A4:H4").Copy()
$wsVal.Range("A6:H6").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Fill in the "Scaled Beta" distribution definition.
$wsVal.Range("A6").Value2 = "Scaled Beta"
$wsVal.Range("B6").Value2 = 4
$wsVal.Range("C6").Value2 = "Min (a)"
$wsVal.Range("D6").Value2 = "Max (b)"
$wsVal.Range("E6").Value2 = "Alpha"
$wsVal.Range("F6").Value2 = "Beta"

$wsVal.Range("F15").Select()

# --- Extend the named ranges that back the validation dropdown/lookup by one row ---
$wb.Names.Item("Validation_Distribution_Parameter_Count").RefersTo = "='Internal - Data Validation'!`$A`$2:`$B`$12"
$wb.Names.Item("Validation_Distribution_Types").RefersTo = "='Internal - Data Validation'!`$A`$2:`$A`$12"

# --- Restore the previously-active sheet/selection ---
$wsExtent = $wb.Worksheets.Item("Extent of Contamination")
$wsExtent.Activate()
$wsExtent.Range("F3").Select()
